# "Generate Report for Handoff"
# The localization-status report is regenerated: the row that used to
# describe 77993f59-...md (status "In Translation") is updated to reflect
# that it has now been handed off for translation ("Ready for handoff" /
# "mt" priority, with a new handoff timestamp), and the report rows are
# re-emitted with that file first (ahead of 15e5b837-...md).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
#   Row 2 <- 77993f59...md   (unchanged status/date, just reordered)
#   Row 3 <- 15e5b837...md   (status/date updated to the new handoff)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "77993f59-45a5-47a3-867b-2c68a3d313b4.md"
$ws.Range("C2").Value = ".md"
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("G2").Value = "2016-08-13 08:13:43"

$ws.Range("A3").Value = "15e5b837-fc20-4136-9319-1e91b6130996.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-13 08:14:33"

# Hyperlink display text follows the new row contents (the link targets
# stay bound to their original relationship ids).
$links = @($ws.Hyperlinks)
$addrB2 = $links[0].Address
$addrB3 = $links[1].Address
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $addrB2, [Type]::Missing, [Type]::Missing, "e2e\77993f59-45a5-47a3-867b-2c68a3d313b4.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $addrB3, [Type]::Missing, [Type]::Missing, "e2e\15e5b837-fc20-4136-9319-1e91b6130996.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
#   Row 2 <- 77993f59...md   (unchanged status/priority/date)
#   Row 3 <- 15e5b837...md   (status -> Ready for handoff, priority -> mt,
#                              new handoff file timestamp)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "77993f59-45a5-47a3-867b-2c68a3d313b4.md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "ht"
$ws.Range("G2").Value = "77993f59-45a5-47a3-867b-2c68a3d313b4.ea7a203b499650941ca5cd2b78ab3adab0b595cf.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-13 08:13:35"

$ws.Range("A3").Value = "15e5b837-fc20-4136-9319-1e91b6130996.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "15e5b837-fc20-4136-9319-1e91b6130996.37cbd955249380062b2d81075a091df82258156d.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-13 08:14:26"

$links = @($ws.Hyperlinks)
$addrA2 = $links[0].Address
$addrA3 = $links[1].Address
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $addrA2, [Type]::Missing, [Type]::Missing, "77993f59-45a5-47a3-867b-2c68a3d313b4.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $addrA3, [Type]::Missing, [Type]::Missing, "15e5b837-fc20-4136-9319-1e91b6130996.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
#   Row 2 <- 77993f59...md   (unchanged status/priority/date)
#   Row 3 <- 15e5b837...md   (status -> Ready for handoff, priority -> mt,
#                              new handoff file timestamp)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "77993f59-45a5-47a3-867b-2c68a3d313b4.md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "ht"
$ws.Range("G2").Value = "77993f59-45a5-47a3-867b-2c68a3d313b4.ea7a203b499650941ca5cd2b78ab3adab0b595cf.de-de.xlf"
$ws.Range("H2").Value = "2016-08-13 08:13:43"

$ws.Range("A3").Value = "15e5b837-fc20-4136-9319-1e91b6130996.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "15e5b837-fc20-4136-9319-1e91b6130996.37cbd955249380062b2d81075a091df82258156d.de-de.xlf"
$ws.Range("H3").Value = "2016-08-13 08:14:33"

$links = @($ws.Hyperlinks)
$addrA2 = $links[0].Address
$addrA3 = $links[1].Address
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $addrA2, [Type]::Missing, [Type]::Missing, "77993f59-45a5-47a3-867b-2c68a3d313b4.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $addrA3, [Type]::Missing, [Type]::Missing, "15e5b837-fc20-4136-9319-1e91b6130996.md")
